$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value2 = 0.005
$ws.Range("F2").Value2 = 0.004
$ws.Range("G2").Value2 = 0.007
$ws.Range("E3").Value2 = 0.006
$ws.Range("F3").Value2 = 0.005
$ws.Range("G3").Value2 = 0.008
$ws.Range("E4").Value2 = 0.007
$ws.Range("G4").Value2 = 0.009
$ws.Range("E5").Value2 = 0.009
$ws.Range("F5").Value2 = 0.007
$ws.Range("G5").Value2 = 0.011
$ws.Range("E6").Value2 = 0.011
$ws.Range("F6").Value2 = 0.009
$ws.Range("G6").Value2 = 0.013
$ws.Range("E7").Value2 = 0.014
$ws.Range("F7").Value2 = 0.011
$ws.Range("E8").Value2 = 0.017
$ws.Range("F8").Value2 = 0.013
$ws.Range("G8").Value2 = 0.019
$ws.Range("E9").Value2 = 0.022
$ws.Range("F9").Value2 = 0.017
$ws.Range("G9").Value2 = 0.025
$ws.Range("E10").Value2 = 0.029
$ws.Range("F10").Value2 = 0.022
$ws.Range("G10").Value2 = 0.031
$ws.Range("E11").Value2 = 0.037
$ws.Range("F11").Value2 = 0.029
$ws.Range("G11").Value2 = 0.041
$ws.Range("E12").Value2 = 0.049
$ws.Range("F12").Value2 = 0.038
$ws.Range("G12").Value2 = 0.054
$ws.Range("E13").Value2 = 0.064
$ws.Range("F13").Value2 = 0.051
$ws.Range("G13").Value2 = 0.072
$ws.Range("E14").Value2 = 0.086
$ws.Range("F14").Value2 = 0.069
$ws.Range("G14").Value2 = 0.096
$ws.Range("E15").Value2 = 0.113
$ws.Range("F15").Value2 = 0.092
$ws.Range("G15").Value2 = 0.126
$ws.Range("E16").Value2 = 0.148
$ws.Range("F16").Value2 = 0.122
$ws.Range("G16").Value2 = 0.166
$ws.Range("E17").Value2 = 0.193
$ws.Range("F17").Value2 = 0.16
$ws.Range("G17").Value2 = 0.208
$ws.Range("E18").Value2 = 0.226
$ws.Range("F18").Value2 = 0.201
$ws.Range("G18").Value2 = 0.25
$ws.Range("E19").Value2 = 0.28
$ws.Range("F19").Value2 = 0.243
$ws.Range("G19").Value2 = 0.312
$ws.Range("E20").Value2 = 0.341
$ws.Range("F20").Value2 = 0.303
$ws.Range("G20").Value2 = 0.376
$ws.Range("E21").Value2 = 0.395
$ws.Range("F21").Value2 = 0.366
$ws.Range("G21").Value2 = 0.431
$ws.Range("E22").Value2 = 0.444
$ws.Range("F22").Value2 = 0.418
$ws.Range("G22").Value2 = 0.483
$ws.Range("E23").Value2 = 0.49
$ws.Range("F23").Value2 = 0.469
$ws.Range("G23").Value2 = 0.544
$ws.Range("E24").Value2 = 0.546
$ws.Range("F24").Value2 = 0.524
$ws.Range("G24").Value2 = 0.602
$ws.Range("E25").Value2 = 0.596
$ws.Range("F25").Value2 = 0.584
$ws.Range("G25").Value2 = 0.657
$ws.Range("E26").Value2 = 0.637
$ws.Range("F26").Value2 = 0.637
$ws.Range("G26").Value2 = 0.712
$ws.Range("E27").Value2 = 0.69
$ws.Range("F27").Value2 = 0.688
$ws.Range("G27").Value2 = 0.778
$ws.Range("E28").Value2 = 0.73
$ws.Range("F28").Value2 = 0.753
$ws.Range("G28").Value2 = 0.816
$ws.Range("E29").Value2 = 0.743
$ws.Range("F29").Value2 = 0.789
$ws.Range("G29").Value2 = 0.841
$ws.Range("E30").Value2 = 0.758
$ws.Range("F30").Value2 = 0.811
$ws.Range("G30").Value2 = 0.866
$ws.Range("E31").Value2 = 0.774
$ws.Range("F31").Value2 = 0.834
$ws.Range("G31").Value2 = 0.891
$ws.Range("E32").Value2 = 0.788
$ws.Range("F32").Value2 = 0.857
$ws.Range("G32").Value2 = 0.916
$ws.Range("E33").Value2 = 0.801
$ws.Range("F33").Value2 = 0.88
$ws.Range("G33").Value2 = 0.94
$ws.Range("E34").Value2 = 0.816
$ws.Range("F34").Value2 = 0.903
$ws.Range("G34").Value2 = 0.964
$ws.Range("E35").Value2 = 0.831
$ws.Range("F35").Value2 = 0.926
$ws.Range("G35").Value2 = 0.992
$ws.Range("E36").Value2 = 0.853
$ws.Range("F36").Value2 = 0.948
$ws.Range("G36").Value2 = 1.015
$ws.Range("E37").Value2 = 0.868
$ws.Range("F37").Value2 = 0.976
$ws.Range("G37").Value2 = 1.035
$ws.Range("E38").Value2 = 0.881
$ws.Range("F38").Value2 = 0.996
$ws.Range("G38").Value2 = 1.053
$ws.Range("E39").Value2 = 0.895
$ws.Range("F39").Value2 = 1.016
$ws.Range("G39").Value2 = 1.071
$ws.Range("E40").Value2 = 0.91
$ws.Range("F40").Value2 = 1.035
$ws.Range("G40").Value2 = 1.089
$ws.Range("E41").Value2 = 0.925
$ws.Range("F41").Value2 = 1.052
$ws.Range("G41").Value2 = 1.108
$ws.Range("E42").Value2 = 0.938
$ws.Range("F42").Value2 = 1.07
$ws.Range("G42").Value2 = 1.123
$ws.Range("E43").Value2 = 0.952
$ws.Range("F43").Value2 = 1.088
$ws.Range("G43").Value2 = 1.14
$ws.Range("E44").Value2 = 0.967
$ws.Range("F44").Value2 = 1.104
$ws.Range("G44").Value2 = 1.155
$ws.Range("E45").Value2 = 0.983
$ws.Range("F45").Value2 = 1.121
$ws.Range("G45").Value2 = 1.171
$ws.Range("E46").Value2 = 0.999
$ws.Range("F46").Value2 = 1.136
$ws.Range("G46").Value2 = 1.188
$ws.Range("E47").Value2 = 1.018
$ws.Range("F47").Value2 = 1.151
$ws.Range("G47").Value2 = 1.203
$ws.Range("E48").Value2 = 1.029
$ws.Range("F48").Value2 = 1.166
$ws.Range("G48").Value2 = 1.221
$ws.Range("E49").Value2 = 1.045
$ws.Range("F49").Value2 = 1.179
$ws.Range("E50").Value2 = 1.062
$ws.Range("F50").Value2 = 1.194
$ws.Range("G50").Value2 = 1.25
$ws.Range("E51").Value2 = 1.081
$ws.Range("F51").Value2 = 1.207
$ws.Range("G51").Value2 = 1.264
$ws.Range("E52").Value2 = 1.091
$ws.Range("F52").Value2 = 1.22
$ws.Range("G52").Value2 = 1.277
$ws.Range("E53").Value2 = 1.106
$ws.Range("F53").Value2 = 1.236
$ws.Range("G53").Value2 = 1.291
$ws.Range("E54").Value2 = 1.118
$ws.Range("F54").Value2 = 1.248
$ws.Range("G54").Value2 = 1.304
$ws.Range("E55").Value2 = 1.131
$ws.Range("F55").Value2 = 1.262
$ws.Range("G55").Value2 = 1.314
$ws.Range("E56").Value2 = 1.142
$ws.Range("F56").Value2 = 1.273
$ws.Range("G56").Value2 = 1.326
$ws.Range("E57").Value2 = 1.153
$ws.Range("G57").Value2 = 1.339
$ws.Range("E58").Value2 = 1.166
$ws.Range("F58").Value2 = 1.294
$ws.Range("G58").Value2 = 1.347
$ws.Range("E59").Value2 = 1.178
$ws.Range("F59").Value2 = 1.306
$ws.Range("G59").Value2 = 1.36
$ws.Range("E60").Value2 = 1.189
$ws.Range("F60").Value2 = 1.316
$ws.Range("G60").Value2 = 1.368
$ws.Range("E61").Value2 = 1.2
$ws.Range("F61").Value2 = 1.326
$ws.Range("G61").Value2 = 1.379
$ws.Range("E62").Value2 = 1.211
$ws.Range("F62").Value2 = 1.336
$ws.Range("G62").Value2 = 1.389
$ws.Range("E63").Value2 = 1.223
$ws.Range("F63").Value2 = 1.347
$ws.Range("G63").Value2 = 1.397
$ws.Range("E64").Value2 = 1.233
$ws.Range("F64").Value2 = 1.355
$ws.Range("G64").Value2 = 1.406
$ws.Range("E65").Value2 = 1.243
$ws.Range("F65").Value2 = 1.364
$ws.Range("G65").Value2 = 1.413
$ws.Range("E66").Value2 = 1.252
$ws.Range("F66").Value2 = 1.372
$ws.Range("G66").Value2 = 1.423

$ws.Range("J9").Select()
